$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 31

# Copy formatting from the row above (row 30) into the new row so the
# new data inherits the same styles used throughout the table.
$ws.Range("A30:J30").Copy() | Out-Null
$ws.Range("A$newRow:J$newRow").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = "DefaultHitTime"
$ws.Cells.Item($newRow, 2).Value = "float"
$ws.Cells.Item($newRow, 3).Value = $false
$ws.Cells.Item($newRow, 4).Value = $false
$ws.Cells.Item($newRow, 5).Value = $false
$ws.Cells.Item($newRow, 6).Value = $true
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = "Friend"
$ws.Cells.Item($newRow, 10).Value = "缺省打击时间（本来应该打到但是物理没碰撞到或者其他原因）"

$ws.Range("J32").Select() | Out-Null
